# "creation of test cases for Student Management"
# Adds new "Sign Ups" and "Updating of Student Details" test-case rows
# (rows 20-26) to the "Iteration 2" sheet, restyles the remaining blank
# rows (27-39) to match, updates one existing cell (E19), widens column D,
# and moves the active selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Widen column D (21.50499916 -> ~29.12999916 width units)
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 28.35

# ---------------------------------------------------------------------
# 2. Fix up the wording of the existing E19 cell
#    ("... List of Tutor Pages." -> "... List of Tutor Page.")
# ---------------------------------------------------------------------
$ws.Range("E19").Value = "Login as an admin user and navigate to the List of Tutor Page. Input the details into the corresponding field and click Create Tutor Select the tutor Amy Tan and press Delete"

# ---------------------------------------------------------------------
# 3. Re-format the whole B20:H39 block to the "bordered / centred /
#    wrapped" look that already decorates the header-style rows
#    (mirrors cell style used on A20:A39 / row 19).
# ---------------------------------------------------------------------
$fmtRange = $ws.Range("B20:H39")
$fmtRange.HorizontalAlignment = -4108
$fmtRange.VerticalAlignment = -4108
$fmtRange.WrapText = $true

# ---------------------------------------------------------------------
# 4. Populate the new "Sign Ups" test cases (rows 20-22)
# ---------------------------------------------------------------------
$ws.Range("B20").Value = "Sign Ups"
$ws.Range("C20").Value = "Verify that when all the neccessary information are filled in, the admin user would be able to successfully create the Student record in the database"
$ws.Range("D20").Value = "Name: Jenny Kim`nID: T0018765H`nAge: 14`nGender: F`nPhone: 91978630`nAddress: Bukit Panjang Ring Rd`nSubjects: Maths, Science`nRequired Amt: `$320`nOutstanding Amt: `$0"
$ws.Range("E20").Value = "Login as an admin user and navigate to the Sign Ups Page. Input the information as specified into the corresponding fields and click Sign Up"
$ws.Range("F20").Value = "Student Creation Successful"

$ws.Range("B21").Value = "Sign Ups"
$ws.Range("C21").Value = "Verify that when mandatory fields are empty, the admin user would not be able to successfully create the Student record in the database"
$ws.Range("D21").Value = "Name: Jenny Kim`nID: `nAge: 14`nGender: F`nPhone: 91978630`nAddress: Bukit Panjang Ring Rd`nSubjects: Maths, Science`nRequired Amt: `$320`nOutstanding Amt: `$0"
$ws.Range("E21").Value = "Login as an admin user and navigate to the Sign Ups Page. Input the information as specified into the corresponding fields and click Sign Up"
$ws.Range("F21").Value = "Error Message regarding Student Creation Failure to be displayed"

$ws.Range("B22").Value = "Sign Ups"
$ws.Range("C22").Value = "Verify that when there is a duplicate of student, the admin user would not be able to overwrite the existing record in the database"
$ws.Range("D22").Value = "Name: Jenny Kim`nID: T0018765H`nAge: 14`nGender: F`nPhone: 91978630`nAddress: Bukit Panjang Ring Rd`nSubjects: Maths, Science`nRequired Amt: `$320`nOutstanding Amt: `$0"
$ws.Range("E22").Value = "Login as an admin user and navigate to the Sign Ups Page. Input the information as specified into the corresponding fields and click Sign Up"
$ws.Range("F22").Value = "Error Message prompting that there was already a record of the specified student"

# ---------------------------------------------------------------------
# 5. Populate the new "Updating of Student Details" test cases (rows 23-26)
# ---------------------------------------------------------------------
$ws.Range("B23").Value = "Updating of Student Details"
$ws.Range("C23").Value = "Verify that if the user is logged in as an admin user, he would be able to modify the records of a specified tutor"
$ws.Range("D23").Value = "-"
$ws.Range("E23").Value = "Login as an admin user and navigate to the Students Records Page. Select the student Jenny Kim and edit her phone number before clicking Update"
$ws.Range("F23").Value = "Update Successful"

$ws.Range("B24").Value = "Updating of Student Details"
$ws.Range("C24").Value = "Verify that if the new  phone number entered does not satisfy the requirements of 8 digits, the update would not be successful "
$ws.Range("D24").Value = "Phone No. : 8333999"
$ws.Range("E24").Value = "Login as an admin user and navigate to the Students Records Page. Select the student Jenny Kim and edit her phone number as specified before clicking Update"
$ws.Range("F24").Value = "Update Unsuccessful"

$ws.Range("B25").Value = "Updating of Student Details"
$ws.Range("C25").Value = "Verify that if the new  ID entered does not satisfy the format requirements, the update would not be successful "
$ws.Range("D25").Value = "ID : T833399"
$ws.Range("E25").Value = "Login as an admin user and navigate to the Students Records Page. Select the student Jenny Kim and edit her ID as specified before clicking Update"
$ws.Range("F25").Value = "Update Unsuccessful"

$ws.Range("B26").Value = "Updating of Student Details"
$ws.Range("C26").Value = "Verify that if the new  subjects entered is not being offered by the tuition centre, the update would not be successful "
$ws.Range("D26").Value = "Subjects : Chinese"
$ws.Range("E26").Value = "Login as an admin user and navigate to the Students Records Page. Select the student Jenny Kim and edit her subjects taken as specified before clicking Update"
$ws.Range("F26").Value = "Update Unsuccessful"

# ---------------------------------------------------------------------
# 6. Row heights: the three "Sign Ups" rows and the four "Updating of
#    Student Details" rows need to grow to fit their wrapped content.
# ---------------------------------------------------------------------
$ws.Rows.Item(20).RowHeight = 147.7
$ws.Rows.Item(21).RowHeight = 147.7
$ws.Rows.Item(22).RowHeight = 147.7
$ws.Rows.Item(23).RowHeight = 115
$ws.Rows.Item(24).RowHeight = 115
$ws.Rows.Item(25).RowHeight = 115
$ws.Rows.Item(26).RowHeight = 115

# ---------------------------------------------------------------------
# 7. Move the active selection / scroll position to show the newly
#    added rows (mirrors what the author was doing at commit time).
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B27:F27").Select() | Out-Null
